# FeedBack Survey Form F9.xlsx - "feat: sops Update 4"
#
# This workbook moved from the "5-Software Service Catalog -SS" folder to the
# "1-Software Development Lifecycle" folder, and the sheet/print-area were
# renamed from the "S-SW-SC-09" (Software Service Catalog) naming scheme to
# the "F-SW-SD-09" (Software Development) naming scheme. The footer revision
# stamp and the current selection/scroll position were also updated.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# 1) Rename the main sheet: S-SW-SC-09 -> F-SW-SD-09
$ws.Name = "F-SW-SD-09"

# 2) Update the workbook-level Print_Area defined name to reference the
#    renamed sheet (renaming the sheet does not retarget defined names here).
foreach ($n in $wb.Names) {
    if ($n.Name -like "*Print_Area*") {
        $n.RefersTo = "='F-SW-SD-09'!`$A`$1:`$F`$29"
    }
}

# 3) Move the selection / view down to the C18:F18 entry (was G19).
$ws.Activate() | Out-Null
$ws.Range("A17").Select() | Out-Null
$excel.ActiveWindow.ScrollRow = 17
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("C18:F18").Select() | Out-Null

# 4) Bump the footer revision stamp: "Rev: 0(0/0/2025)" -> "Rev:0(01/10/2025)"
$ws.PageSetup.RightFooter = '&"Arial,Regular"&16Rev:0(01/10/2025)'

$wb.Save() | Out-Null
